$wb = $excel.ActiveWorkbook

# Sheet "展览" - update 想去人数 (F column) values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 845
$wsExhibit.Range("F5").Value = 1053
$wsExhibit.Range("F6").Value = 2446
$wsExhibit.Range("F7").Value = 207

# Sheet "全部类型" - update 想去人数 (F column) values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 845
$wsAll.Range("F7").Value = 1053
$wsAll.Range("F8").Value = 2446
$wsAll.Range("F10").Value = 207
